# Re-run of the HP-toolbox Lacex fit after changing the Tarrival lower-bound
# threshold (T5: 30 -> 20). The underlying MATLAB simulation was re-run with
# the new bound, producing refreshed fit results in A1:M16; this script
# replays the resulting cell values (the MATLAB recompute itself happened
# out-of-band and isn't reproducible from within Excel).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated lower bound for Tarrival (row 5 holds the lower bounds used by
# the R-column "Lower/Ok/Upper" checks for data column E) ---
$ws.Range("T5").Value = 20

# --- Refreshed simulation/fit output values (A1:M16) ---
$newValues = @(
    @("A1", "1.0000000002246367E-4"),
    @("B1", "1.3285975171052873E-2"),
    @("C1", "4.5410137928182014E-2"),
    @("D1", "150000189.85702357"),
    @("E1", "27.864607276962733"),
    @("F1", "54.780848970299019"),
    @("G1", "1.0000031066132442E-8"),
    @("H1", "0.59997658802644682"),
    @("I1", "3.5252117786944458E-5"),
    @("J1", "2.0593439124861734E-3"),
    @("K1", "0.97747298300119134"),
    @("L1", "-0.49344324702422182"),
    @("M1", "-0.7000409581160838"),
    @("A2", "8.6661738012672254E-3"),
    @("B2", "0.55138682912517523"),
    @("C2", "2.9037123200575041E-2"),
    @("D2", "119779666.15910518"),
    @("E2", "23.725269400855069"),
    @("F2", "66.797791756301621"),
    @("G2", "0.20082226663264988"),
    @("H2", "9.8559033304461089E-3"),
    @("I2", "5.9333725380924401E-4"),
    @("J2", "2.0300946260441272E-3"),
    @("K2", "0.99574120770731744"),
    @("L2", "0.98234265166820556"),
    @("M2", "0.85119438142528137"),
    @("A3", "1.0000000002225244E-4"),
    @("B3", "1.1276816868432057E-2"),
    @("C3", "4.5445811351612325E-2"),
    @("D3", "150000064.9195846"),
    @("E3", "31.724804762290479"),
    @("F3", "54.138341213074924"),
    @("G3", "1.000002319736726E-8"),
    @("H3", "0.59999767516501779"),
    @("I3", "4.2762327658428232E-6"),
    @("J3", "2.1777410517568661E-3"),
    @("K3", "0.80577342477386971"),
    @("L3", "-0.53941290306549528"),
    @("M3", "-0.73654476567460114"),
    @("A4", "1.0000000002230694E-4"),
    @("B4", "1.1923457734909134E-2"),
    @("C4", "4.5434921242829877E-2"),
    @("D4", "150000239.15305385"),
    @("E4", "40.223257240810938"),
    @("F4", "58.920625547013493"),
    @("G4", "1.0000024565277651E-8"),
    @("H4", "0.59999191750123682"),
    @("I4", "8.4393671752232924E-6"),
    @("J4", "2.2245635329467745E-3"),
    @("K4", "0.94052103268490395"),
    @("L4", "-0.5135081540415023"),
    @("M4", "-0.55863801367132693"),
    @("A5", "1.0000000002223794E-4"),
    @("B5", "1.2953341672327022E-2"),
    @("C5", "4.5454232198170089E-2"),
    @("D5", "150000061.27472425"),
    @("E5", "30.13531567047368"),
    @("F5", "53.099336295663889"),
    @("G5", "1.0000022878916901E-8"),
    @("H5", "0.59999727561987581"),
    @("I5", "1.1164378939968163E-6"),
    @("J5", "2.1670145505032884E-3"),
    @("K5", "0.72742231695528181"),
    @("L5", "-0.50812670802918358"),
    @("M5", "-0.57010968486757618"),
    @("A6", "1.0000000002232024E-4"),
    @("B6", "1.2754881983591192E-2"),
    @("C6", "4.5426718009964652E-2"),
    @("D6", "150000126.61099541"),
    @("E6", "33.849702615348249"),
    @("F6", "54.057395301640149"),
    @("G6", "1.0000024947115184E-8"),
    @("H6", "0.59998958714933626"),
    @("I6", "3.6852456601258068E-6"),
    @("J6", "2.1771367654063742E-3"),
    @("K6", "0.90983020278869997"),
    @("L6", "-0.55509703523853116"),
    @("M6", "-0.60925714656751584"),
    @("A7", "1.0838045914027154E-2"),
    @("B7", "6.4573083837011621E-3"),
    @("C7", "3.6328217467279614E-2"),
    @("D7", "150167675.68179601"),
    @("E7", "39.349767877357422"),
    @("F7", "53.622138874369853"),
    @("G7", "1.0000023371612956E-8"),
    @("H7", "0.59093617805999521"),
    @("I7", "4.4722197271424057E-4"),
    @("J7", "1.4842614835726567E-3"),
    @("K7", "0.99519687803839696"),
    @("L7", "0.37065768000301202"),
    @("M7", "-0.4973145661254601"),
    @("A8", "1.0000000002237555E-4"),
    @("B8", "9.977035289958993E-3"),
    @("C8", "4.5396596646422066E-2"),
    @("D8", "150000110.42974392"),
    @("E8", "38.052157675133742"),
    @("F8", "55.883937936202393"),
    @("G8", "1.0000026799409753E-8"),
    @("H8", "0.59999573435467635"),
    @("I8", "1.8028253437368054E-5"),
    @("J8", "2.2098698077807366E-3"),
    @("K8", "0.84451555542505219"),
    @("L8", "-0.49791529130878764"),
    @("M8", "-0.54857722322415126"),
    @("A9", "1.000000000222197E-4"),
    @("B9", "9.9988362889332147E-3"),
    @("C9", "4.5451423246041685E-2"),
    @("D9", "150000066.6048401"),
    @("E9", "39.698431270358313"),
    @("F9", "57.076034214667722"),
    @("G9", "1.0000022501492048E-8"),
    @("H9", "0.59999965956140111"),
    @("I9", "1.1366127428115871E-6"),
    @("J9", "2.2461344756480044E-3"),
    @("K9", "0.26439263722395356"),
    @("L9", "-0.35804828516822207"),
    @("M9", "-0.44984156059141367"),
    @("A10", "7.5375029476333838E-3"),
    @("B10", "1.2160449333588226"),
    @("C10", "2.6395590591137646E-2"),
    @("D10", "100000000.19232847"),
    @("E10", "40.308698710908992"),
    @("F10", "45.861777643223306"),
    @("G10", "1.4771705413739111E-3"),
    @("H10", "2.8504565668281689E-4"),
    @("I10", "0.89867699189212102"),
    @("J10", "9.9994127085341236E-2"),
    @("K10", "0.99116143819587532"),
    @("L10", "0.9927631432494205"),
    @("M10", "0.96119607482708536"),
    @("A11", "2.6209723748675851E-2"),
    @("B11", "0.52934684405812649"),
    @("C11", "2.6594247555100314E-2"),
    @("D11", "150142161.52967486"),
    @("E11", "20.001489305243467"),
    @("F11", "71.114467047765714"),
    @("G11", "1.734657631652692"),
    @("H11", "8.5423133097486362E-4"),
    @("I11", "0.81815536921747956"),
    @("J11", "4.3704411573062795E-5"),
    @("K11", "0.98366025995805773"),
    @("L11", "0.95713040570866215"),
    @("M11", "0.88037277833742988"),
    @("A12", "1.0000000002228612E-4"),
    @("B12", "9.990825750683453E-3"),
    @("C12", "4.5428033194248714E-2"),
    @("D12", "150000052.14751464"),
    @("E12", "33.522864423273056"),
    @("F12", "55.62635282915371"),
    @("G12", "1.0000024007651981E-8"),
    @("H12", "0.59999760960633708"),
    @("I12", "6.7542007776368318E-6"),
    @("J12", "2.2066719671406535E-3"),
    @("K12", "0.34276506061467971"),
    @("L12", "-0.56161187672943491"),
    @("M12", "-0.67242145913056861"),
    @("A13", "1.0000000002225389E-4"),
    @("B13", "9.9920626871110967E-3"),
    @("C13", "4.5436405686289788E-2"),
    @("D13", "150000045.57887048"),
    @("E13", "38.51605222815131"),
    @("F13", "55.940495420089682"),
    @("G13", "1.0000023230178751E-8"),
    @("H13", "0.59999920660193795"),
    @("I13", "6.170544489694359E-6"),
    @("J13", "2.2331315235188075E-3"),
    @("K13", "0.32960756924186119"),
    @("L13", "-0.54316690902499243"),
    @("M13", "-0.60270975135112992"),
    @("A14", "1.0000000002226187E-4"),
    @("B14", "9.9899014004585958E-3"),
    @("C14", "4.5431861195288276E-2"),
    @("D14", "150000058.63665292"),
    @("E14", "37.376388175555164"),
    @("F14", "55.115057197678439"),
    @("G14", "1.000002341418302E-8"),
    @("H14", "0.59999908789323131"),
    @("I14", "7.4512182206618473E-6"),
    @("J14", "2.2245317278717274E-3"),
    @("K14", "0.31097626356553654"),
    @("L14", "-0.58019565790975869"),
    @("M14", "-0.6178308567205304"),
    @("A15", "1.0000000002225956E-4"),
    @("B15", "9.9928864554518874E-3"),
    @("C15", "4.5434609385371853E-2"),
    @("D15", "150000066.6664634"),
    @("E15", "37.656203632759429"),
    @("F15", "56.175773159118464"),
    @("G15", "1.0000023360201194E-8"),
    @("H15", "0.59999829022618767"),
    @("I15", "5.4986632366613223E-6"),
    @("J15", "2.229600533894139E-3"),
    @("K15", "0.30353902707807179"),
    @("L15", "-0.64191603743893677"),
    @("M15", "-0.69907881706239472"),
    @("A16", "1.0000000002222338E-4"),
    @("B16", "1.0172533382061075E-2"),
    @("C16", "4.5451845117286858E-2"),
    @("D16", "150000074.01884153"),
    @("E16", "43.264765786440861"),
    @("F16", "59.369875287587526"),
    @("G16", "1.0000022576001589E-8"),
    @("H16", "0.59999965856787663"),
    @("I16", "1.4269196932959207E-6"),
    @("J16", "2.2662815923935853E-3"),
    @("K16", "0.56852519020501635"),
    @("L16", "-0.36182202498591076"),
    @("M16", "-0.47512930455471092")
)

foreach ($pair in $newValues) {
    $cellRef = $pair[0]
    $cellVal = $pair[1]
    $ws.Range($cellRef).Value = [double]$cellVal
}

# --- Restore the viewport/selection as left by the author (scrolled back to
# column A, active cell Q10) ---
$ws.Range("Q10").Select()

Write-Host "Applied Tarrival threshold change and refreshed simulation outputs."
